# "Overlapping times are correctly calculated" [Fixes #1]
#
# Two participant rows had the wrong student ID / e-mail copied into
# them, which threw off the overlap calculation between sessions:
#   - Row 8  ("강강강")  had someone else's ID and e-mail -> fix to her own
#                        ID/e-mail (the same ones already used in row 4).
#   - Row 9  ("황황황")  had the wrong e-mail address      -> fix to the
#                        e-mail that row 8 used to (incorrectly) have.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the mismatched ID number / e-mail -------------------------
$ws.Range("A8").Value = "2021-12344"
$ws.Range("C8").Value = "test2@snu.ac.kr"
$ws.Range("C9").Value = "test6@snu.ac.kr"

# --- Re-point the column-C e-mail hyperlinks to match ------------------
# A plain value overwrite does not retarget a cell's existing hyperlink,
# so rebuild the full set (the engine only supports clearing every
# hyperlink on the sheet at once) in the same order Excel keeps them.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:test@naver.com")   | Out-Null
$ws.Range("C4").Copy()
$ws.Range("C2").PasteSpecial(-4122)   # xlPasteFormats - keep the Hyperlink cell style intact

$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:test2@snu.ac.kr")  | Out-Null
$ws.Range("C5").Copy()
$ws.Range("C4").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:test3@snu.ac.kr")  | Out-Null
$ws.Range("C6").Copy()
$ws.Range("C5").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:test4@snu.ac.kr")  | Out-Null
$ws.Range("C7").Copy()
$ws.Range("C6").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:test5@snu.ac.kr")  | Out-Null
$ws.Range("C9").Copy()
$ws.Range("C7").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:test6@snu.ac.kr")  | Out-Null
$ws.Range("C10").Copy()
$ws.Range("C9").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:test7@snu.ac.kr") | Out-Null
$ws.Range("C8").Copy()
$ws.Range("C10").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:test2@snu.ac.kr")  | Out-Null
$ws.Range("C2").Copy()
$ws.Range("C8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Widen the Join time / Leave time columns so the (now meaningful)
#     overlap windows are easy to read -----------------------------------
$ws.Columns.Item(4).ColumnWidth = 18.498697916666668
$ws.Columns.Item(5).ColumnWidth = 20.166666666666668

# --- Leave the selection on the corrected cell --------------------------
$ws.Range("C9").Select()
